$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Add new "CreatorUserID" column (F) with header and values
$ws.Range("F1").Value = "CreatorUserID"
$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 5
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 4
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = 3
$ws.Range("F8").Value = 5
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 1

# Match the column width used for the new column
$ws.Columns.Item(6).ColumnWidth = 19.17

# Make the Issues sheet the active tab with F1 selected
$ws.Activate() | Out-Null
$ws.Range("F1").Select() | Out-Null
